# Weekly refresh of the fruit/vegetable price rows: the 20 data rows
# (rows 2-21) get their Fecha/Volumen/Precio columns reshuffled among
# themselves (row 20 happens to land back on itself).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (read source's old values, write into destination)
$rowMap = @{
    2  = 17
    3  = 21
    4  = 10
    5  = 9
    6  = 15
    7  = 6
    8  = 4
    9  = 5
    10 = 3
    11 = 13
    12 = 16
    13 = 2
    14 = 19
    15 = 7
    16 = 14
    17 = 8
    18 = 12
    19 = 11
    20 = 20
    21 = 18
}

$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Snapshot the old values for every touched column/row before writing anything,
# since destinations and sources overlap.
$old = @{}
foreach ($r in $rowMap.Keys) {
    foreach ($c in $cols) {
        $old["$r,$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $old["$srcRow,$c"]
    }
}
